$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert rows 15-17 with new data (Data / Quantidade de horas)
# Use raw date serial numbers so no time-of-day fraction is introduced
$ws.Range("A15").Value = 41561
$ws.Range("B15").Value = 2.5

$ws.Range("A16").Value = 41562
$ws.Range("B16").Value = 1

$ws.Range("A17").Value = 41563
$ws.Range("B17").Value = 5

# Apply the date style (numFmtId 14 - short date) matching other date cells in column A
$ws.Range("A2").Copy()
$ws.Range("A15:A17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Update selection to A18
$ws.Range("A18").Select()

$wb.Save()
